$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The numeric placeholder ids in column A (1..6) get replaced with the real
# customer names, and a 7th customer (AVK) is appended as a new row.
# Assign in the order the strings were first introduced so the regenerated
# shared-strings table lines up with the source file (A2, then the brand-new
# row 8, then the remaining rows bottom-to-top before A3).
$ws.Range("A2").Value = "Vetafarm"
$ws.Range("A8").Value = "AVK"
$ws.Range("A5").Value = "Vulcan"
$ws.Range("A6").Value = "AU Ramps"
$ws.Range("A7").Value = "TTI"
$ws.Range("A4").Value = "Commodore"
$ws.Range("A3").Value = "3T Foods"

# Fill in the rest of the new row (row 8 / AVK) - length, depth, height, weight.
$ws.Range("B8").Value = 1.2
$ws.Range("C8").Value = 1.2
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 600

# Column A needs to be widened (and best-fit) to comfortably show the longest
# name now that it holds text instead of single-digit ids.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(1).ColumnWidth = 10.33

# The last user selection recorded in the saved file was cell A4.
$ws.Range("A4").Select() | Out-Null
